# Apply updated res_line/pl_mw.xlsx values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.324950471376837
$ws.Range("C2").Value = 0.2719280110609645
$ws.Range("D2").Value = 0.07906340444534976
$ws.Range("E2").Value = 0.09033757861834069
$ws.Range("G2").Value = 0.4196179103708673
$ws.Range("H2").Value = 0.5737809257208824
$ws.Range("I2").Value = 0.5111246775245561
$ws.Range("L2").Value = 0.2244772498705601
$ws.Range("O2").Value = 1.935940658437289

# Row 3
$ws.Range("B3").Value = 1.181201200701878
$ws.Range("C3").Value = 0.2531973612183833
$ws.Range("D3").Value = 0.07164508350165022
$ws.Range("E3").Value = 0.09155319993447364
$ws.Range("G3").Value = 0.4213809614588158
$ws.Range("H3").Value = 0.579350118437354
$ws.Range("I3").Value = 0.5210108412055323
$ws.Range("L3").Value = 0.2137058924912054
$ws.Range("O3").Value = 1.951024393143385

# Row 4
$ws.Range("B4").Value = 1.092754359451249
$ws.Range("C4").Value = 0.2416516024323414
$ws.Range("D4").Value = 0.06712429838881917
$ws.Range("E4").Value = 0.09235810668841449
$ws.Range("G4").Value = 0.4229137285292595
$ws.Range("H4").Value = 0.5831390204242268
$ws.Range("I4").Value = 0.5275213011869848
$ws.Range("L4").Value = 0.2071904858615028
$ws.Range("O4").Value = 1.96200277369779

# Row 5
$ws.Range("B5").Value = 1.056667614873561
$ws.Range("C5").Value = 0.2369356416293158
$ws.Range("D5").Value = 0.06529061278416748
$ws.Range("E5").Value = 0.0927008213505518
$ws.Range("G5").Value = 0.4236512021034926
$ws.Range("H5").Value = 0.5847758471420832
$ws.Range("I5").Value = 0.5302848664483637
$ws.Range("L5").Value = 0.2045602147786383
$ws.Range("O5").Value = 1.966907368903378

# Row 6
$ws.Range("B6").Value = 1.050672855176742
$ws.Range("C6").Value = 0.2361519080492087
$ws.Range("D6").Value = 0.0649866482996373
$ws.Range("E6").Value = 0.09275861709757205
$ws.Range("G6").Value = 0.4237804648435173
$ws.Range("H6").Value = 0.5850532449501458
$ws.Range("I6").Value = 0.5307504217921331
$ws.Range("L6").Value = 0.2041249618228562
$ws.Range("O6").Value = 1.96774776658124

# Row 7
$ws.Range("B7").Value = 1.092267855726732
$ws.Range("C7").Value = 0.2415880452724082
$ws.Range("D7").Value = 0.06709953396297408
$ws.Range("E7").Value = 0.09236266910918722
$ws.Range("G7").Value = 0.4229232178744979
$ws.Range("H7").Value = 0.5831607195166697
$ws.Range("I7").Value = 0.5275581244787091
$ws.Range("L7").Value = 0.2071549125114842
$ws.Range("O7").Value = 1.962067175644492

# Row 8
$ws.Range("B8").Value = 1.275425333931139
$ws.Range("C8").Value = 0.2654792461468389
$ws.Range("D8").Value = 0.07649848758399003
$ws.Range("E8").Value = 0.0907445785391463
$ws.Range("G8").Value = 0.4201321447527135
$ws.Range("H8").Value = 0.575624473841664
$ws.Range("I8").Value = 0.5144419075881821
$ws.Range("L8").Value = 0.2207429524939073
$ws.Range("O8").Value = 1.940784634901306

# Row 9
$ws.Range("B9").Value = 1.633048350162255
$ws.Range("C9").Value = 0.3119586989034246
$ws.Range("D9").Value = 0.09520141555778139
$ws.Range("E9").Value = 0.08803593980819002
$ws.Range("G9").Value = 0.4182480823735659
$ws.Range("H9").Value = 0.5637801636573272
$ws.Range("I9").Value = 0.4922243609679064
$ws.Range("L9").Value = 0.248165644051511
$ws.Range("O9").Value = 1.912715118293193

# Row 10
$ws.Range("B10").Value = 1.894760208384298
$ws.Range("C10").Value = 0.3458660996918184
$ws.Range("D10").Value = 0.1091109700846857
$ws.Range("E10").Value = 0.08632926031089561
$ws.Range("G10").Value = 0.419075805411893
$ws.Range("H10").Value = 0.5568716706959265
$ws.Range("I10").Value = 0.4780496291047029
$ws.Range("L10").Value = 0.2687848604384442
$ws.Range("O10").Value = 1.900484222891663

# Row 11
$ws.Range("B11").Value = 2.013578107629087
$ws.Range("C11").Value = 0.3612362324419394
$ws.Range("D11").Value = 0.1154760850311192
$ws.Range("E11").Value = 0.08561442460214685
$ws.Range("G11").Value = 0.4199378054465086
$ws.Range("H11").Value = 0.5541193010551524
$ws.Range("I11").Value = 0.4720704227970813
$ws.Range("L11").Value = 0.2782673178865309
$ws.Range("O11").Value = 1.896755441858545

# Row 12
$ws.Range("B12").Value = 2.058535414670644
$ws.Range("C12").Value = 0.3670483534741322
$ws.Range("D12").Value = 0.1178918171712553
$ws.Range("E12").Value = 0.08535258968889892
$ws.Range("G12").Value = 0.420334430161077
$ws.Range("H12").Value = 0.5531332646048526
$ws.Range("I12").Value = 0.4698739236217868
$ws.Range("L12").Value = 0.28187277366861
$ws.Range("O12").Value = 1.895608365091363

# Row 13
$ws.Range("B13").Value = 2.048854714955723
$ws.Range("C13").Value = 0.3657969805068717
$ws.Range("D13").Value = 0.1173713058787058
$ws.Range("E13").Value = 0.085408586435749
$ws.Range("G13").Value = 0.4202458813275456
$ws.Range("H13").Value = 0.553343122996111
$ws.Range("I13").Value = 0.4703439640354539
$ws.Range("L13").Value = 0.2810956242081346
$ws.Range("O13").Value = 1.895843609116412

# Row 14
$ws.Range("B14").Value = 2.01727751712076
$ws.Range("C14").Value = 0.361714565909125
$ws.Range("D14").Value = 0.1156747205436375
$ws.Range("E14").Value = 0.08559270571442745
$ws.Range("G14").Value = 0.4199690267543872
$ws.Range("H14").Value = 0.5540370518374544
$ws.Range("I14").Value = 0.471888357251359
$ws.Range("L14").Value = 0.2785636478210449
$ws.Range("O14").Value = 1.896655755574216

# Row 15
$ws.Range("B15").Value = 1.997930750758314
$ws.Range("C15").Value = 0.359212888216689
$ws.Range("D15").Value = 0.1146362163133574
$ws.Range("E15").Value = 0.08570663788776756
$ws.Range("G15").Value = 0.4198085996464869
$ws.Range("H15").Value = 0.5544694284433831
$ws.Range("I15").Value = 0.4728431665344281
$ws.Range("L15").Value = 0.2770146444471067
$ws.Range("O15").Value = 1.897187751558192

# Row 16
$ws.Range("B16").Value = 1.886990213609636
$ws.Range("C16").Value = 0.3448604962118225
$ws.Range("D16").Value = 0.1086957509095754
$ws.Range("E16").Value = 0.08637721558489986
$ws.Range("G16").Value = 0.4190292742711961
$ws.Range("H16").Value = 0.5570594073473387
$ws.Range("I16").Value = 0.4784498431906776
$ws.Range("L16").Value = 0.2681672162274111
$ws.Range("O16").Value = 1.900764924107278

# Row 17
$ws.Range("B17").Value = 1.818869523284718
$ws.Range("C17").Value = 0.336041537769745
$ws.Range("D17").Value = 0.1050610925402964
$ws.Range("E17").Value = 0.08680436141748693
$ws.Range("G17").Value = 0.4186758041549865
$ws.Range("H17").Value = 0.5587483189795535
$ws.Range("I17").Value = 0.4820096537894116
$ws.Range("L17").Value = 0.2627658245895361
$ws.Range("O17").Value = 1.903430165390176

# Row 18
$ws.Range("B18").Value = 1.779666225275037
$ws.Range("C18").Value = 0.3309639972244156
$ws.Range("D18").Value = 0.102974071492028
$ws.Range("E18").Value = 0.08705583647636139
$ws.Range("G18").Value = 0.4185181682671839
$ws.Range("H18").Value = 0.5597564710627267
$ws.Range("I18").Value = 0.4841012868931571
$ws.Range("L18").Value = 0.2596687575041585
$ws.Range("O18").Value = 1.905135809260855

# Row 19
$ws.Range("B19").Value = 1.766388946591121
$ws.Range("C19").Value = 0.3292439632902244
$ws.Range("D19").Value = 0.1022680498454207
$ws.Range("E19").Value = 0.08714197615094754
$ws.Range("G19").Value = 0.418472627398927
$ws.Range("H19").Value = 0.560104120340327
$ws.Range("I19").Value = 0.4848170483041727
$ws.Range("L19").Value = 0.2586218080614771
$ws.Range("O19").Value = 1.905742934883932

# Row 20
$ws.Range("B20").Value = 1.826123390276905
$ws.Range("C20").Value = 0.3369808624367465
$ws.Range("D20").Value = 0.1054476419540435
$ws.Range("E20").Value = 0.08675829149634673
$ws.Range("G20").Value = 0.4187087020311395
$ws.Range("H20").Value = 0.5585647288429527
$ws.Range("I20").Value = 0.4816261372834489
$ws.Range("L20").Value = 0.2633398118816501
$ws.Range("O20").Value = 1.90312856787196

# Row 21
$ws.Range("B21").Value = 2.026553514663931
$ws.Range("C21").Value = 0.3629138959539944
$ws.Range("D21").Value = 0.1161729022290388
$ws.Range("E21").Value = 0.08553838495392085
$ws.Range("G21").Value = 0.4200484372016291
$ws.Range("H21").Value = 0.5538317013894272
$ws.Range("I21").Value = 0.4714328921793367
$ws.Range("L21").Value = 0.2793069537989936
$ws.Range("O21").Value = 1.896410010006576

# Row 22
$ws.Range("B22").Value = 2.157332570910114
$ws.Range("C22").Value = 0.3798145795481105
$ws.Range("D22").Value = 0.1232139602090427
$ws.Range("E22").Value = 0.08479273966102241
$ws.Range("G22").Value = 0.4213334158991273
$ws.Range("H22").Value = 0.5510661679652742
$ws.Range("I22").Value = 0.4651657252861163
$ws.Range("L22").Value = 0.2898277739428465
$ws.Range("O22").Value = 1.893563749458025

# Row 23
$ws.Range("B23").Value = 2.087553692113033
$ws.Range("C23").Value = 0.3707988866038363
$ws.Range("D23").Value = 0.1194531360506659
$ws.Range("E23").Value = 0.08518597734071065
$ws.Range("G23").Value = 0.4206100100335846
$ws.Range("H23").Value = 0.5525121622627296
$ws.Range("I23").Value = 0.4684744280067363
$ws.Range("L23").Value = 0.2842048388459517
$ws.Range("O23").Value = 1.89494116797357

# Row 24
$ws.Range("B24").Value = 1.822844040812868
$ws.Range("C24").Value = 0.3365562167900578
$ws.Range("D24").Value = 0.1052728748895646
$ws.Range("E24").Value = 0.08677910132038136
$ws.Range("G24").Value = 0.4186936869563169
$ws.Range("H24").Value = 0.5586476141536707
$ws.Range("I24").Value = 0.4817993847844875
$ws.Range("L24").Value = 0.2630802864763524
$ws.Range("O24").Value = 1.903264380176438

# Row 25
$ws.Range("B25").Value = 1.536476986031971
$ws.Range("C25").Value = 0.2994260522797276
$ws.Range("D25").Value = 0.09011239538018856
$ws.Range("E25").Value = 0.08871895397963847
$ws.Range("G25").Value = 0.418371046290261
$ws.Range("H25").Value = 0.5666697426199505
$ws.Range("I25").Value = 0.4978584183444106
$ws.Range("L25").Value = 0.240664114622362
$ws.Range("O25").Value = 1.918839364982034

Write-Host "Updated $($wb.Name): 216 cells across rows 2-25 (columns B,C,D,E,G,H,I,L,O)"
